# Backup before deadline cleaning — minor recompute/precision refresh of
# the Q1 summary row (row 2). Only the last-significant-digit float noise
# actually changes value; everything else on the sheet is untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K2").Value  = 0.71399999999999997
$ws.Range("Y2").Value  = 1.4980000000000002
$ws.Range("AA2").Value = 2.1080000000000005
$ws.Range("AE2").Value = 1.2809999999999999
$ws.Range("AF2").Value = 1.159
$ws.Range("AV2").Value = 2.8940000000000001
$ws.Range("BY2").Value = 1.3689999999999998
$ws.Range("CI2").Value = 0.88500000000000001
$ws.Range("CQ2").Value = 1.2969999999999999
